$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N2").Value = "Actual"
$ws.Range("N3").Value = "# Slides"
$ws.Range("N4").Formula = "=4+4"
$ws.Range("N6").Value = 5
$ws.Range("N8").Value = 7
$ws.Range("I3").ClearFormats()
$ws.Range("K3").ClearFormats()
$ws.Range("L3").ClearFormats()
$ws.Range("C21").ClearFormats()
$ws.Range("N8").Select()
